$wb = $excel.ActiveWorkbook

# --- settings sheet: bump form_version and add showFooter = 1 ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20210221001
$settings.Range("A9").Value = "showFooter"
$settings.Range("B9").Value = 1

# --- make settings the active sheet / tab, with B10 selected ---
$settings.Activate()
$settings.Range("B10").Select()

$wb.Save()
